$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Fix spelling mistake: "Fäll" -> "fällt"
$ws.Range("D2").Value = "fällt"

# Update selection (reflects new active cell after edits / print-area review)
$ws.Range("E12").Select()
